$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Value = "HelloTest"
Write-Host "value2:" $ws.Range("A1").Value2
